$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览": update "想去人数" (F column) values ---
$ws1.Range("F2").Value = 72
$ws1.Range("F3").Value = 243
$ws1.Range("F4").Value = 508
$ws1.Range("F5").Value = 2217
$ws1.Range("F7").Value = 8149
$ws1.Range("F8").Value = 99
$ws1.Range("F10").Value = 267
$ws1.Range("F11").Value = 1793
$ws1.Range("F12").Value = 1588
$ws1.Range("F14").Value = 205
$ws1.Range("F15").Value = 4341
$ws1.Range("F16").Value = 6156
$ws1.Range("F17").Value = 757
$ws1.Range("F19").Value = 1188
$ws1.Range("F20").Value = 1257
$ws1.Range("F21").Value = 464
$ws1.Range("F22").Value = 6405
$ws1.Range("F23").Value = 356
$ws1.Range("F25").Value = 4341
$ws1.Range("F26").Value = 295
$ws1.Range("F27").Value = 713
$ws1.Range("F28").Value = 2005
$ws1.Range("F29").Value = 1187
$ws1.Range("F30").Value = 333
$ws1.Range("F31").Value = 1070
$ws1.Range("F32").Value = 45
$ws1.Range("F33").Value = 43
$ws1.Range("F34").Value = 76
$ws1.Range("F36").Value = 1183
$ws1.Range("F37").Value = 1897
$ws1.Range("F38").Value = 130
$ws1.Range("F39").Value = 433
$ws1.Range("F40").Value = 167
$ws1.Range("F41").Value = 1198
$ws1.Range("F42").Value = 563
$ws1.Range("F44").Value = 1136
$ws1.Range("F47").Value = 185
$ws1.Range("F48").Value = 26
$ws1.Range("F49").Value = 10

# --- Sheet "演出": remove 4 rows (old rows 4-7) that were deleted, shifting rows up ---
$ws2.Range("A4:I7").Delete(-4162)

# --- Sheet "演出": fix the sequential index column (A) after the shift ---
for ($r = 4; $r -le 37; $r++) {
    $ws2.Range("A$r").Value = $r - 1
}

# --- Sheet "演出": update "想去人数" (F column) values (using post-shift row numbers) ---
$ws2.Range("F7").Value = 20
$ws2.Range("F9").Value = 690
$ws2.Range("F10").Value = 390
$ws2.Range("F11").Value = 410
$ws2.Range("F14").Value = 113
$ws2.Range("F18").Value = 185
$ws2.Range("F20").Value = 104
$ws2.Range("F23").Value = 104
$ws2.Range("F25").Value = 129
$ws2.Range("F28").Value = 281
$ws2.Range("F29").Value = 100
$ws2.Range("F34").Value = 10

# --- Sheet "本地生活": update "想去人数" (F column) values ---
$ws3.Range("F6").Value = 1588
$ws3.Range("F7").Value = 489
$ws3.Range("F8").Value = 3141
$ws3.Range("F9").Value = 1034
$ws3.Range("F10").Value = 1137
$ws3.Range("F11").Value = 1413
$ws3.Range("F12").Value = 1778
$ws3.Range("F13").Value = 258
$ws3.Range("F14").Value = 111

# --- Sheet "全部类型": update "想去人数" (F column) values ---
$ws4.Range("F3").Value = 243
$ws4.Range("F5").Value = 508
$ws4.Range("F6").Value = 489
$ws4.Range("F7").Value = 3141
$ws4.Range("F8").Value = 2217
$ws4.Range("F9").Value = 99
$ws4.Range("F10").Value = 1034
$ws4.Range("F12").Value = 20
$ws4.Range("F13").Value = 267
$ws4.Range("F14").Value = 1793
$ws4.Range("F15").Value = 1588
$ws4.Range("F16").Value = 1413
$ws4.Range("F17").Value = 690
$ws4.Range("F18").Value = 205
$ws4.Range("F19").Value = 1778
$ws4.Range("F20").Value = 4342
$ws4.Range("F21").Value = 390
$ws4.Range("F22").Value = 410
$ws4.Range("F23").Value = 757
$ws4.Range("F25").Value = 1188
$ws4.Range("F26").Value = 1257
$ws4.Range("F27").Value = 464
$ws4.Range("F28").Value = 6405
$ws4.Range("F29").Value = 356
$ws4.Range("F30").Value = 295
$ws4.Range("F31").Value = 713
$ws4.Range("F32").Value = 2005
$ws4.Range("F33").Value = 1187
$ws4.Range("F34").Value = 333
$ws4.Range("F35").Value = 43
$ws4.Range("F36").Value = 76
$ws4.Range("F38").Value = 1897
$ws4.Range("F39").Value = 130
$ws4.Range("F40").Value = 433
$ws4.Range("F41").Value = 104
$ws4.Range("F42").Value = 167
$ws4.Range("F43").Value = 1198
$ws4.Range("F44").Value = 129
$ws4.Range("F45").Value = 563
$ws4.Range("F46").Value = 282
$ws4.Range("F47").Value = 1136
$ws4.Range("F48").Value = 185
$ws4.Range("F49").Value = 26
